$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("B2").Value = 0.1519756838905775
$ws.Range("C2").Value = 0.6048632218844985
$ws.Range("J2").Value = 0.0182370820668693
$ws.Range("O2").Value = 0.00303951367781155
$ws.Range("P2").Value = 0.1306990881458966
$ws.Range("S2").Value = 0.0911854103343465
# Row 3
$ws.Range("B3").Value = 0.004524886877828055
$ws.Range("C3").Value = 0.03167420814479638
$ws.Range("J3").Value = 0.03619909502262444
$ws.Range("P3").Value = 0.7194570135746606
$ws.Range("S3").Value = 0.2081447963800905
# Row 4
$ws.Range("J4").Value = 0.07142857142857142
$ws.Range("P4").Value = 0.5714285714285714
$ws.Range("S4").Value = 0.3571428571428572
# Row 5
$ws.Range("J5").Value = 0.3333333333333333
$ws.Range("P5").Value = 0.6666666666666666
# Row 6
$ws.Range("B6").Value = 0.06734006734006734
$ws.Range("D6").Value = 0.0101010101010101
$ws.Range("E6").Value = 0.003367003367003367
$ws.Range("F6").Value = 0.07744107744107744
$ws.Range("J6").Value = 0.2457912457912458
$ws.Range("O6").Value = 0.03367003367003367
$ws.Range("Q6").Value = 0.1077441077441077
$ws.Range("R6").Value = 0.08080808080808081
$ws.Range("S6").Value = 0.3737373737373738
# Row 7
$ws.Range("B7").Value = 0.1203319502074689
$ws.Range("D7").Value = 0.04149377593360996
$ws.Range("F7").Value = 0.05394190871369295
$ws.Range("J7").Value = 0.1078838174273859
$ws.Range("O7").Value = 0.04149377593360996
$ws.Range("Q7").Value = 0.1327800829875519
$ws.Range("R7").Value = 0.1037344398340249
$ws.Range("S7").Value = 0.3983402489626556
# Row 8
$ws.Range("B8").Value = 0.124031007751938
$ws.Range("D8").Value = 0.007751937984496124
$ws.Range("E8").Value = 0.001937984496124031
$ws.Range("F8").Value = 0.04457364341085272
$ws.Range("J8").Value = 0.1298449612403101
$ws.Range("O8").Value = 0.01937984496124031
$ws.Range("Q8").Value = 0.1182170542635659
$ws.Range("R8").Value = 0.1046511627906977
$ws.Range("S8").Value = 0.4496124031007752
# Row 9
$ws.Range("B9").Value = 0.09848484848484848
$ws.Range("D9").Value = 0.02272727272727273
$ws.Range("E9").Value = 0.007575757575757576
$ws.Range("F9").Value = 0.09090909090909091
$ws.Range("J9").Value = 0.1060606060606061
$ws.Range("O9").Value = 0.01515151515151515
$ws.Range("Q9").Value = 0.1287878787878788
$ws.Range("R9").Value = 0.09848484848484848
$ws.Range("S9").Value = 0.4318181818181818
# Row 10
$ws.Range("B10").Value = 0.1053811659192825
$ws.Range("D10").Value = 0.0179372197309417
$ws.Range("F10").Value = 0.09491778774289986
$ws.Range("J10").Value = 0.1210762331838565
$ws.Range("O10").Value = 0.01943198804185351
$ws.Range("Q10").Value = 0.187593423019432
$ws.Range("R10").Value = 0.08520179372197309
$ws.Range("S10").Value = 0.3684603886397608
# Row 11
$ws.Range("G11").Value = 0.1945205479452055
$ws.Range("J11").Value = 0.08493150684931507
$ws.Range("K11").Value = 0.2438356164383562
$ws.Range("L11").Value = 0.4684931506849315
$ws.Range("S11").Value = 0.00821917808219178
# Row 12
$ws.Range("G12").Value = 0.7151162790697675
$ws.Range("J12").Value = 0.2325581395348837
$ws.Range("K12").Value = 0.01162790697674419
$ws.Range("S12").Value = 0.04069767441860465
# Row 13
$ws.Range("G13").Value = 0.6753246753246753
$ws.Range("J13").Value = 0.2987012987012987
$ws.Range("S13").Value = 0.02597402597402598
# Row 15
$ws.Range("F15").Value = 0.04147465437788019
$ws.Range("H15").Value = 0.1566820276497696
$ws.Range("I15").Value = 0.04147465437788019
$ws.Range("J15").Value = 0.3271889400921659
$ws.Range("K15").Value = 0.06912442396313365
$ws.Range("M15").Value = 0.009216589861751152
$ws.Range("O15").Value = 0.05069124423963134
$ws.Range("S15").Value = 0.304147465437788
# Row 16
$ws.Range("F16").Value = 0.02727272727272727
$ws.Range("H16").Value = 0.1727272727272727
$ws.Range("I16").Value = 0.09545454545454546
$ws.Range("J16").Value = 0.3772727272727273
$ws.Range("K16").Value = 0.1045454545454545
$ws.Range("M16").Value = 0.05
$ws.Range("O16").Value = 0.02727272727272727
$ws.Range("S16").Value = 0.1454545454545454
# Row 17
$ws.Range("F17").Value = 0.02255639097744361
$ws.Range("H17").Value = 0.1604010025062657
$ws.Range("I17").Value = 0.04761904761904762
$ws.Range("J17").Value = 0.4411027568922306
$ws.Range("K17").Value = 0.09523809523809523
$ws.Range("M17").Value = 0.02506265664160401
$ws.Range("O17").Value = 0.07268170426065163
$ws.Range("S17").Value = 0.1353383458646616
# Row 18
$ws.Range("F18").Value = 0.02586206896551724
$ws.Range("H18").Value = 0.2155172413793103
$ws.Range("I18").Value = 0.05603448275862069
$ws.Range("J18").Value = 0.418103448275862
$ws.Range("K18").Value = 0.1163793103448276
$ws.Range("M18").Value = 0.01724137931034483
$ws.Range("O18").Value = 0.04310344827586207
$ws.Range("S18").Value = 0.1077586206896552
# Row 19
$ws.Range("F19").Value = 0.02723463687150838
$ws.Range("H19").Value = 0.2325418994413408
$ws.Range("I19").Value = 0.04818435754189944
$ws.Range("J19").Value = 0.3456703910614525
$ws.Range("K19").Value = 0.1201117318435754
$ws.Range("M19").Value = 0.03631284916201118
$ws.Range("N19").Value = 0.001396648044692737
$ws.Range("O19").Value = 0.05726256983240224
$ws.Range("S19").Value = 0.1312849162011173

$wb.Save()
